# Removed Test Case Inter-Dependency
$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name on the input sheet (B1) and short name (B2)
$wsInput.Range("B1").Value = "4255-MS-EI-DB-DL-REC-CTRFD-RNI-INT-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PE-1st"
$wsInput.Range("B2").Value = "425f"

# Keep output sheet in sync with the new product name
$wsOutput.Range("B1").Value = "4255-MS-EI-DB-DL-REC-CTRFD-RNI-INT-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PE-1st"

# Move the active selection on the input sheet from B17 to B2
$wsInput.Range("B2").Select()

# Make the output sheet the active tab, removing the dependency on
# ProductLoanInput's selection state
$wsOutput.Activate()
